$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43 - this shifts existing rows 43..151 down to 44..152,
# matching the target diff (a new weekly price record was added at the top of this block).
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record's data.
$ws.Cells.Item(43, 1).Value2  = 6
$ws.Cells.Item(43, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(43, 3).Value2  = "Metropolitana"
$ws.Cells.Item(43, 4).Value2  = 45281
$ws.Cells.Item(43, 5).Value2  = 13
$ws.Cells.Item(43, 6).Value2  = 100114007
$ws.Cells.Item(43, 7).Value2  = "Jengibre"
$ws.Cells.Item(43, 8).Value2  = "Sin especificar"
$ws.Cells.Item(43, 9).Value2  = "Primera"
$ws.Cells.Item(43, 10).Value2 = 200
$ws.Cells.Item(43, 11).Value2 = 25000
$ws.Cells.Item(43, 12).Value2 = 25000
$ws.Cells.Item(43, 13).Value2 = 25000
$ws.Cells.Item(43, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(43, 15).Value2 = "Perú"
$ws.Cells.Item(43, 16).Value2 = 1923
$ws.Cells.Item(43, 17).Value2 = 13
$ws.Cells.Item(43, 18).Value2 = "Hortaliza"
